# Stackup updated with recommendations from fab engineer.

$wb = $excel.ActiveWorkbook

# Remove the empty "Sheet3" worksheet entirely.
$wb.Worksheets.Item("Sheet3").Delete() | Out-Null

# Work on the "8-layer" sheet.
$ws = $wb.Worksheets.Item("8-layer")

# Narrow column E a bit per fab recommendations.
$ws.Columns.Item(5).ColumnWidth = 6.5

# New trace-width / calculated-Z0 entries for rows 2 and 16 (outer layers).
$ws.Cells.Item(2, 8).Value = 5
$ws.Cells.Item(2, 9).Value = 88

# Updated trace width / calculated Z0 for rows 4 and 14.
$ws.Cells.Item(4, 8).Value = 5
$ws.Cells.Item(4, 9).Value = 49.9

$ws.Cells.Item(14, 8).Value = 5
$ws.Cells.Item(14, 9).Value = 49.9

$ws.Cells.Item(16, 8).Value = 5
$ws.Cells.Item(16, 9).Value = 88

# Update the saved selection to match the author's final cursor position.
$ws.Range("G19").Select() | Out-Null
